$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Strava module "outputs" cell (E2) to include the new
# "date range" line, reflecting the updated output of the strava module.
$ws.Range("E2").Value = "number of entries for leaderboard for all segments per time of day`ncoordinates of segments`ndate range"
